# chore(runtime): publish files + archive (2025-12-16 15:07:52)
#
# The source data was refreshed: the KHL game played on 01-12-2025
# (Torpedo vs Dinamo M, game_uid 897837) was removed from the
# "Form_Games" sheet, which shifts all subsequent rows up by one and
# shrinks the used range from A1:AH81 to A1:AH80. The "Aggregates"
# sheet's per-team summary row for "Torpedo" is then refreshed to
# reflect the updated set of games for that team.

$wb = $excel.ActiveWorkbook

# --- 1. Remove the stale game row from Form_Games -----------------------
$wsGames = $wb.Worksheets.Item("Form_Games")
$wsGames.Rows.Item(74).Delete()

# --- 2. Refresh the recalculated Torpedo aggregate row -------------------
$wsAgg = $wb.Worksheets.Item("Aggregates")

$row = 16
$wsAgg.Cells.Item($row, 2).Value  = 4                     # rows
$wsAgg.Cells.Item($row, 3).Value  = 3                      # GF_mean
$wsAgg.Cells.Item($row, 4).Value  = 2.75                   # GA_mean
$wsAgg.Cells.Item($row, 5).Value  = 0.2                    # GF_trend
$wsAgg.Cells.Item($row, 6).Value  = 0.3                    # GA_trend
$wsAgg.Cells.Item($row, 7).Value  = 30                     # SOG_for_mean
$wsAgg.Cells.Item($row, 8).Value  = 80.66666666666667      # SOG_for_var
$wsAgg.Cells.Item($row, 9).Value  = -0.4                   # SOG_for_trend
$wsAgg.Cells.Item($row, 10).Value = 29                     # SOG_against_mean
$wsAgg.Cells.Item($row, 11).Value = 14                     # SOG_against_var
$wsAgg.Cells.Item($row, 12).Value = 1.4                    # SOG_against_trend
$wsAgg.Cells.Item($row, 13).Value = 20.25                  # HITS_mean
$wsAgg.Cells.Item($row, 14).Value = 2.9                    # HITS_trend
$wsAgg.Cells.Item($row, 15).Value = 49.6                   # FOW_pct_mean
$wsAgg.Cells.Item($row, 16).Value = -4.320000000000001     # FOW_pct_trend
$wsAgg.Cells.Item($row, 17).Value = 0.2075                 # PP_eff_by_minutes_mean
$wsAgg.Cells.Item($row, 18).Value = -0.049                 # PP_eff_by_minutes_trend
$wsAgg.Cells.Item($row, 19).Value = 0.5825                 # PK_eff_by_minutes_mean
$wsAgg.Cells.Item($row, 20).Value = 0.101                  # PK_eff_by_minutes_trend
$wsAgg.Cells.Item($row, 21).Value = 2                       # INT_mean
$wsAgg.Cells.Item($row, 22).Value = -0.4                   # INT_trend
$wsAgg.Cells.Item($row, 23).Value = 6                       # PIM_mean
$wsAgg.Cells.Item($row, 24).Value = 1.6                    # PIM_trend
$wsAgg.Cells.Item($row, 25).Value = 9                       # PIM_opp_mean
$wsAgg.Cells.Item($row, 26).Value = 0.4                    # PIM_opp_trend
$wsAgg.Cells.Item($row, 27).Value = 10.5                   # SH_pct_mean
$wsAgg.Cells.Item($row, 28).Value = 0.78                   # SH_pct_trend
